$wb = $excel.ActiveWorkbook

# --- Sheet 1: pvERC_15y ---
$ws1 = $wb.Worksheets.Item("pvERC_15y")
$ws1.Range("L2").Value = -13.425011624025345
$ws1.Range("L3").Value = -12.444300712312028
$ws1.Range("L4").Value = -10.901462620174312
$ws1.Range("L5").Value = -7.646489363993291
$ws1.Range("L6").Value = -1.0531870573520952

# --- Sheet 2: pvERC_30y ---
$ws2 = $wb.Worksheets.Item("pvERC_30y")
$ws2.Range("L2").Value = -20.15713359428066
$ws2.Range("L3").Value = -17.226754581133356
$ws2.Range("L4").Value = -13.54864984128766
$ws2.Range("L5").Value = -10.33776641864193
$ws2.Range("L6").Value = -4.459986107743418
